# Generate Report for Handback
#
# The en-US source file (487e3f7f-bd11-43fc-a557-eb6bdf7faa6c.md) has now been
# handed back for both the zh-cn and de-de locales: its status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns are now populated for that row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$targetFileName   = "487e3f7f-bd11-43fc-a557-eb6bdf7faa6c.md"
$targetFileUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a3e8117d5e1309781292b05c443b421ec2d6b3a/e2e/487e3f7f-bd11-43fc-a557-eb6bdf7faa6c.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetFileUrl, "", "", $targetFileName)
$wsZh.Range("J2").Value = $wsZh.Range("G2").Value2
$wsZh.Range("K2").Value = "2016-08-27 08:39:49"

$wsZh.Columns.Item(3).ColumnWidth = 29.09
$wsZh.Columns.Item(9).ColumnWidth = 38.92
$wsZh.Columns.Item(10).ColumnWidth = 39.09

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetFileUrl, "", "", $targetFileName)
$wsDe.Range("J2").Value = $wsDe.Range("G2").Value2
$wsDe.Range("K2").Value = "2016-08-27 08:39:55"

$wsDe.Columns.Item(3).ColumnWidth = 29.09
$wsDe.Columns.Item(9).ColumnWidth = 38.92
$wsDe.Columns.Item(10).ColumnWidth = 39.09

# ---------------------------------------------------------------------------
# Overview sheet - only the zh-cn / de-de status columns widen (autofit)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09
